$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.177257299423218
$ws.Range("B1").Value = 2.414568662643433
$ws.Range("D1").Value = 2.33542537689209
$ws.Range("E1").Value = 1.20007848739624
